# Reproduces the "Updated symbol list" crypto data refresh (15 -> 16 snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + its new literal text. Coin-name/link cells (B, C) are
# plain text already; price/volume/hour cells (D, E, G) look numeric to Excel, so we
# force the cell to Text (NumberFormat "@") before writing, then drop the leftover
# number-format override with ClearFormats() so no stray style sticks around.
$updates = @(
    @{ Cell = 'D2'; Value = '304.02' }
    @{ Cell = 'E2'; Value = '0.26%' }
    @{ Cell = 'G2'; Value = '16' }
    @{ Cell = 'D3'; Value = '37.13' }
    @{ Cell = 'E3'; Value = '3.24%' }
    @{ Cell = 'G3'; Value = '16' }
    @{ Cell = 'D4'; Value = '5.035' }
    @{ Cell = 'E4'; Value = '-1.10%' }
    @{ Cell = 'G4'; Value = '16' }
    @{ Cell = 'D5'; Value = '0.07816' }
    @{ Cell = 'E5'; Value = '-0.48%' }
    @{ Cell = 'G5'; Value = '16' }
    @{ Cell = 'D6'; Value = '2.249' }
    @{ Cell = 'E6'; Value = '-1.60%' }
    @{ Cell = 'G6'; Value = '16' }
    @{ Cell = 'D7'; Value = '7.985' }
    @{ Cell = 'E7'; Value = '-0.98%' }
    @{ Cell = 'G7'; Value = '16' }
    @{ Cell = 'B8'; Value = 'MXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D8'; Value = '0.9247' }
    @{ Cell = 'E8'; Value = '0.03%' }
    @{ Cell = 'G8'; Value = '16' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D9'; Value = '0.09865' }
    @{ Cell = 'E9'; Value = '-2.76%' }
    @{ Cell = 'G9'; Value = '16' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D10'; Value = '0.1880' }
    @{ Cell = 'E10'; Value = '3.01%' }
    @{ Cell = 'G10'; Value = '16' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D11'; Value = '0.08617' }
    @{ Cell = 'E11'; Value = '0.65%' }
    @{ Cell = 'G11'; Value = '16' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D12'; Value = '0.03652' }
    @{ Cell = 'E12'; Value = '7.56%' }
    @{ Cell = 'G12'; Value = '16' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D13'; Value = '0.09948' }
    @{ Cell = 'E13'; Value = '0.56%' }
    @{ Cell = 'G13'; Value = '16' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D14'; Value = '0.001487' }
    @{ Cell = 'E14'; Value = '0.77%' }
    @{ Cell = 'G14'; Value = '16' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D15'; Value = '0.005716' }
    @{ Cell = 'E15'; Value = '2.16%' }
    @{ Cell = 'G15'; Value = '16' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D16'; Value = '3.455' }
    @{ Cell = 'E16'; Value = '-0.83%' }
    @{ Cell = 'G16'; Value = '16' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D17'; Value = '4.033' }
    @{ Cell = 'E17'; Value = '1.03%' }
    @{ Cell = 'G17'; Value = '16' }
    @{ Cell = 'E18'; Value = '17.72%' }
    @{ Cell = 'G18'; Value = '16' }
    @{ Cell = 'D19'; Value = '0.3432' }
    @{ Cell = 'E19'; Value = '0.01%' }
    @{ Cell = 'G19'; Value = '16' }
    @{ Cell = 'D20'; Value = '0.1303' }
    @{ Cell = 'E20'; Value = '-1.78%' }
    @{ Cell = 'G20'; Value = '16' }
    @{ Cell = 'D21'; Value = '4.775' }
    @{ Cell = 'E21'; Value = '4.89%' }
    @{ Cell = 'G21'; Value = '16' }
    @{ Cell = 'D22'; Value = '0.2206' }
    @{ Cell = 'E22'; Value = '-0.59%' }
    @{ Cell = 'G22'; Value = '16' }
    @{ Cell = 'D23'; Value = '0.04605' }
    @{ Cell = 'E23'; Value = '-0.82%' }
    @{ Cell = 'G23'; Value = '16' }
    @{ Cell = 'D24'; Value = '0.001255' }
    @{ Cell = 'E24'; Value = '1.87%' }
    @{ Cell = 'G24'; Value = '16' }
    @{ Cell = 'D25'; Value = '0.005038' }
    @{ Cell = 'E25'; Value = '12.36%' }
    @{ Cell = 'G25'; Value = '16' }
    @{ Cell = 'D26'; Value = '0.0001405' }
    @{ Cell = 'E26'; Value = '8.14%' }
    @{ Cell = 'G26'; Value = '16' }
    @{ Cell = 'D27'; Value = '0.0002726' }
    @{ Cell = 'E27'; Value = '-19.85%' }
    @{ Cell = 'G27'; Value = '16' }
    @{ Cell = 'G28'; Value = '16' }
    @{ Cell = 'G29'; Value = '16' }
    @{ Cell = 'G30'; Value = '16' }
    @{ Cell = 'G31'; Value = '16' }
    @{ Cell = 'G32'; Value = '16' }
    @{ Cell = 'G33'; Value = '16' }
    @{ Cell = 'G34'; Value = '16' }
    @{ Cell = 'G35'; Value = '16' }
    @{ Cell = 'G36'; Value = '16' }
    @{ Cell = 'G37'; Value = '16' }
    @{ Cell = 'G38'; Value = '16' }
    @{ Cell = 'D39'; Value = '0.01825' }
    @{ Cell = 'E39'; Value = '3.99%' }
    @{ Cell = 'G39'; Value = '16' }
    @{ Cell = 'D40'; Value = '0.04755' }
    @{ Cell = 'E40'; Value = '0.42%' }
    @{ Cell = 'G40'; Value = '16' }
    @{ Cell = 'D41'; Value = '0.007969' }
    @{ Cell = 'E41'; Value = '1.26%' }
    @{ Cell = 'G41'; Value = '16' }
    @{ Cell = 'D42'; Value = '0.1401' }
    @{ Cell = 'E42'; Value = '-0.71%' }
    @{ Cell = 'G42'; Value = '16' }
    @{ Cell = 'D43'; Value = '0.007584' }
    @{ Cell = 'E43'; Value = '-13.85%' }
    @{ Cell = 'G43'; Value = '16' }
    @{ Cell = 'D44'; Value = '0.002151' }
    @{ Cell = 'E44'; Value = '-2.65%' }
    @{ Cell = 'G44'; Value = '16' }
    @{ Cell = 'D45'; Value = '0.01043' }
    @{ Cell = 'E45'; Value = '13.70%' }
    @{ Cell = 'G45'; Value = '16' }
    @{ Cell = 'D46'; Value = '0.00006300' }
    @{ Cell = 'E46'; Value = '5.74%' }
    @{ Cell = 'G46'; Value = '16' }
    @{ Cell = 'E47'; Value = '0.45%' }
    @{ Cell = 'G47'; Value = '16' }
    @{ Cell = 'D48'; Value = '0.0005818' }
    @{ Cell = 'E48'; Value = '0.30%' }
    @{ Cell = 'G48'; Value = '16' }
    @{ Cell = 'D49'; Value = '35.45' }
    @{ Cell = 'E49'; Value = '510.85%' }
    @{ Cell = 'G49'; Value = '16' }
    @{ Cell = 'D50'; Value = '0.002698' }
    @{ Cell = 'E50'; Value = '0.43%' }
    @{ Cell = 'G50'; Value = '16' }
    @{ Cell = 'D51'; Value = '0.00002107' }
    @{ Cell = 'E51'; Value = '0.45%' }
    @{ Cell = 'G51'; Value = '16' }
)

foreach ($u in $updates) {
    $col = $u.Cell -replace '[0-9]+$', ''
    $range = $ws.Range($u.Cell)
    if ($col -eq 'B' -or $col -eq 'C') {
        $range.Value = $u.Value
    } else {
        $range.NumberFormat = '@'
        $range.Value = $u.Value
        $range.ClearFormats()
    }
}
